# Remove redundant "Why This Solution?" slide from the presales presentation.
# This streamlines the deck from 10 slides down to 9: Solution Overview already
# covers technical capabilities, Our Partnership Advantage covers "why choose us",
# and Success Stories supplies proof points - so slide 4 ("Why This Solution?")
# is dropped and every following slide simply shifts up one position.

$p = $ppt.ActivePresentation

# Slide 4 (1-based) is "Why This Solution?" - delete it outright.
$p.Slides.Item(4).Delete()
